# Auto-generated Excel COM-interop script
# Applies cached numeric-value corrections to the Excalibur Profits workbook
# (Leve profit/price columns H, I, J, K, L, M, N) across multiple job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value2 = 0
$ws.Range("J87").Value2 = 0
$ws.Range("L87").Value2 = 0
$ws.Range("H90").Value2 = 0
$ws.Range("J90").Value2 = 0
$ws.Range("L90").Value2 = 0
$ws.Range("H132").Value2 = 141942.7
$ws.Range("I132").Value2 = 169359.28
$ws.Range("K132").Value2 = 508077.84
$ws.Range("M132").Value2 = -505547.84
$ws.Range("H133").Value2 = 88999.5
$ws.Range("J133").Value2 = 88999.5
$ws.Range("L133").Value2 = 88999.5
$ws.Range("N133").Value2 = -99119.5
$ws.Range("H136").Value2 = 78921.75
$ws.Range("J136").Value2 = 78921.75
$ws.Range("L136").Value2 = 78921.75
$ws.Range("N136").Value2 = -89121.75
$ws.Range("N87").ClearContents()
$ws.Range("N90").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 5377972.5
$ws.Range("I32").Value2 = 6579617.5
$ws.Range("J32").Value2 = 5911
$ws.Range("K32").Value2 = 6579617.5
$ws.Range("L32").Value2 = 5911
$ws.Range("M32").Value2 = -6579330.5
$ws.Range("N32").Value2 = -6485
$ws.Range("H61").Value2 = 2497.6785
$ws.Range("I61").Value2 = 2283.9048
$ws.Range("K61").Value2 = 2283.9048
$ws.Range("M61").Value2 = -2071.9048
$ws.Range("H62").Value2 = 100249
$ws.Range("J62").Value2 = 100249
$ws.Range("L62").Value2 = 100249
$ws.Range("N62").Value2 = -101497
$ws.Range("H65").Value2 = 100249
$ws.Range("J65").Value2 = 100249
$ws.Range("L65").Value2 = 300747
$ws.Range("N65").Value2 = -306987
$ws.Range("H132").Value2 = 2214.6897
$ws.Range("I132").Value2 = 1966.3673
$ws.Range("K132").Value2 = 5899.1019
$ws.Range("M132").Value2 = -3369.1019
$ws.Range("H136").Value2 = 2497.6785
$ws.Range("I136").Value2 = 2283.9048
$ws.Range("K136").Value2 = 6851.714399999999
$ws.Range("M136").Value2 = -4301.714399999999
$ws.Range("H139").Value2 = 89494.5
$ws.Range("J139").Value2 = 89494.5
$ws.Range("L139").Value2 = 89494.5
$ws.Range("N139").Value2 = -99774.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value2 = 34249.5
$ws.Range("J58").Value2 = 39999
$ws.Range("L58").Value2 = 39999
$ws.Range("N58").Value2 = -40587
$ws.Range("H94").Value2 = 2458
$ws.Range("I94").Value2 = 2099
$ws.Range("J94").Value2 = 2601.6
$ws.Range("K94").Value2 = 2099
$ws.Range("L94").Value2 = 2601.6
$ws.Range("M94").Value2 = -1648
$ws.Range("N94").Value2 = -3503.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 93542.664
$ws.Range("I31").Value2 = 154892.67
$ws.Range("K31").Value2 = 154892.67
$ws.Range("M31").Value2 = -154597.67
$ws.Range("H34").Value2 = 93542.664
$ws.Range("I34").Value2 = 154892.67
$ws.Range("K34").Value2 = 154892.67
$ws.Range("M34").Value2 = -154690.67
$ws.Range("H58").Value2 = 1744.8572
$ws.Range("I58").Value2 = 1061.5555
$ws.Range("K58").Value2 = 1061.5555
$ws.Range("M58").Value2 = -858.5554999999999
$ws.Range("H68").Value2 = 87357.75
$ws.Range("J68").Value2 = 87357.75
$ws.Range("L68").Value2 = 87357.75
$ws.Range("N68").Value2 = -88855.75
$ws.Range("H70").Value2 = 60000
$ws.Range("J70").Value2 = 60000
$ws.Range("L70").Value2 = 60000
$ws.Range("N70").Value2 = -60630
$ws.Range("H71").Value2 = 87357.75
$ws.Range("J71").Value2 = 87357.75
$ws.Range("L71").Value2 = 262073.25
$ws.Range("N71").Value2 = -269561.25
$ws.Range("H73").Value2 = 60000
$ws.Range("J73").Value2 = 60000
$ws.Range("L73").Value2 = 60000
$ws.Range("N73").Value2 = -62184
$ws.Range("H99").Value2 = 3585.2856
$ws.Range("I99").Value2 = 2999.5
$ws.Range("J99").Value2 = 3819.6
$ws.Range("K99").Value2 = 2999.5
$ws.Range("L99").Value2 = 3819.6
$ws.Range("M99").Value2 = -1501.5
$ws.Range("N99").Value2 = -6815.6
$ws.Range("H124").Value2 = 47405.75
$ws.Range("J124").Value2 = 47405.75
$ws.Range("L124").Value2 = 47405.75
$ws.Range("N124").Value2 = -52315.75
$ws.Range("H126").Value2 = 3585.2856
$ws.Range("I126").Value2 = 2999.5
$ws.Range("J126").Value2 = 3819.6
$ws.Range("K126").Value2 = 8998.5
$ws.Range("L126").Value2 = 11458.8
$ws.Range("M126").Value2 = -6528.5
$ws.Range("N126").Value2 = -16398.8
$ws.Range("H132").Value2 = 35715864
$ws.Range("I132").Value2 = 1842.3334
$ws.Range("K132").Value2 = 5527.0002
$ws.Range("M132").Value2 = -2997.0002
$ws.Range("H134").Value2 = 51570.617
$ws.Range("I134").Value2 = 55717.332
$ws.Range("K134").Value2 = 167151.996
$ws.Range("M134").Value2 = -164616.996
$ws.Range("H136").Value2 = 1744.8572
$ws.Range("I136").Value2 = 1061.5555
$ws.Range("K136").Value2 = 3184.6665
$ws.Range("M136").Value2 = -634.6664999999998
$ws.Range("H138").Value2 = 0
$ws.Range("J138").Value2 = 0
$ws.Range("L138").Value2 = 0
$ws.Range("H139").Value2 = 89718.375
$ws.Range("J139").Value2 = 89549.39999999999
$ws.Range("L139").Value2 = 89549.39999999999
$ws.Range("N139").Value2 = -99829.39999999999
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value2 = 98535.08
$ws.Range("J37").Value2 = 98535.08
$ws.Range("L37").Value2 = 295605.24
$ws.Range("N37").Value2 = -295829.24
$ws.Range("H117").Value2 = 5790.6665
$ws.Range("J117").Value2 = 6606.5557
$ws.Range("L117").Value2 = 19819.6671
$ws.Range("N117").Value2 = -26703.6671
$ws.Range("H131").Value2 = 1474.7142
$ws.Range("J131").Value2 = 1543.9546
$ws.Range("L131").Value2 = 4631.8638
$ws.Range("N131").Value2 = -14711.8638
$ws.Range("H132").Value2 = 5033.3335
$ws.Range("J132").Value2 = 5444.4443
$ws.Range("L132").Value2 = 48999.9987
$ws.Range("N132").Value2 = -54059.9987

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value2 = 47638076
$ws.Range("I132").Value2 = 71444300
$ws.Range("J132").Value2 = 25610
$ws.Range("K132").Value2 = 214332900
$ws.Range("L132").Value2 = 76830
$ws.Range("M132").Value2 = -214330370
$ws.Range("N132").Value2 = -81890

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value2 = 2772.5925
$ws.Range("I132").Value2 = 2656.9565
$ws.Range("J132").Value2 = 3437.5
$ws.Range("K132").Value2 = 7970.869499999999
$ws.Range("L132").Value2 = 10312.5
$ws.Range("M132").Value2 = -5440.869499999999
$ws.Range("N132").Value2 = -15372.5
$ws.Range("H136").Value2 = 76633.82000000001
$ws.Range("I136").Value2 = 2475.6
$ws.Range("J136").Value2 = 182574.14
$ws.Range("K136").Value2 = 7426.799999999999
$ws.Range("L136").Value2 = 547722.42
$ws.Range("M136").Value2 = -4876.799999999999
$ws.Range("N136").Value2 = -552822.42

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value2 = 2351.1667
$ws.Range("I122").Value2 = 2057.2222
$ws.Range("J122").Value2 = 3233
$ws.Range("K122").Value2 = 6171.6666
$ws.Range("L122").Value2 = 9699
$ws.Range("M122").Value2 = -3721.6666
$ws.Range("N122").Value2 = -14599
